# Apply automatic tracker update:
#  - row 236: fill in result (G236="Fallo", H236=-1)
#  - append 4 new match rows (278-281)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the resolved result for row 236 ---
$ws.Range("G236").Value = "Fallo"
$ws.Range("H236").Value = -1

# --- Append new tracker rows ---
# columns: row, event_id, fecha, jugador_A, jugador_B, pronostico, cuota
$newRows = @(
    @(278, 14428722, "2025-08-18", "Vit Kopriva", "Sebastian Korda", "Gana Vit Kopriva", 2.75),
    @(279, 14466777, "2025-08-18", "Johannus Monday", "Harold Mayot", "Gana Johannus Monday", 2.5),
    @(280, 14466765, "2025-08-18", "Francesco Maestrelli", "James McCabe", "Gana Francesco Maestrelli", 2.25),
    @(281, 14466848, "2025-08-18", "Carlota Martinez Cirez", "Arantxa Rus", "Gana Carlota Martinez Cirez", 3.5)
)

foreach ($r in $newRows) {
    $row = $r[0]

    $ws.Cells.Item($row, 1).Value = $r[1]

    # Column B ("fecha") holds a date-shaped string, not a real date value.
    # Force text formatting while writing it, then drop the style again so
    # the cell ends up unstyled like every other row, exactly like the
    # original workbook.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 2).ClearFormats()

    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
    $ws.Cells.Item($row, 5).Value = $r[5]
    $ws.Cells.Item($row, 6).Value = $r[6]

    # Columns G ("resultado") and H ("profit") are left blank for new,
    # still-unresolved matches, matching rows 2-277 before they get a result.
}
